$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Ly86"
$ws.Cells.Item(2,3).Value = "Cd180"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(2,7).Value = [double]"0.641009"
$ws.Cells.Item(2,8).Value = [double]"1.923027"
$ws.Cells.Item(2,9).Value = [double]"0.005069072203501787"
$ws.Cells.Item(2,10).Value = [double]"0.005069072203501787"
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = [double]"3.413185"
$ws.Cells.Item(2,14).Value = [double]"10.239555"
$ws.Cells.Item(2,15).Value = [double]"0.01281891839503857"
$ws.Cells.Item(2,16).Value = [double]"0.01281891839503857"
$ws.Cells.Item(2,17).Value = [double]"2.187882303665"
$ws.Cells.Item(2,18).Value = [double]"19.690940732985"
$ws.Cells.Item(2,19).Value = [double]"6.498002291524777e-05"
$ws.Cells.Item(2,20).Value = [double]"6.498002291524776e-05"

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Ly86"
$ws.Cells.Item(3,3).Value = "Cd180"
$ws.Cells.Item(3,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(3,7).Value = [double]"0.641009"
$ws.Cells.Item(3,8).Value = [double]"1.923027"
$ws.Cells.Item(3,9).Value = [double]"0.005069072203501787"
$ws.Cells.Item(3,10).Value = [double]"0.005069072203501787"
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = [double]"71.73578866666666"
$ws.Cells.Item(3,14).Value = [double]"215.207366"
$ws.Cells.Item(3,15).Value = [double]"0.2694185111330716"
$ws.Cells.Item(3,16).Value = [double]"0.2694185111330716"
$ws.Cells.Item(3,17).Value = [double]"45.98328615743134"
$ws.Cells.Item(3,18).Value = [double]"413.849575416882"
$ws.Cells.Item(3,19).Value = [double]"0.00136570188589349"
$ws.Cells.Item(3,20).Value = [double]"0.00136570188589349"

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Ly86"
$ws.Cells.Item(4,3).Value = "Cd180"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(4,7).Value = [double]"0.641009"
$ws.Cells.Item(4,8).Value = [double]"1.923027"
$ws.Cells.Item(4,9).Value = [double]"0.005069072203501787"
$ws.Cells.Item(4,10).Value = [double]"0.005069072203501787"
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = [double]"0.5338986666666666"
$ws.Cells.Item(4,14).Value = [double]"1.601696"
$ws.Cells.Item(4,15).Value = [double]"0.002005166271157263"
$ws.Cells.Item(4,16).Value = [double]"0.002005166271157263"
$ws.Cells.Item(4,17).Value = [double]"0.3422338504213334"
$ws.Cells.Item(4,18).Value = [double]"3.080104653792"
$ws.Cells.Item(4,19).Value = [double]"1.016433260852261e-05"
$ws.Cells.Item(4,20).Value = [double]"1.016433260852261e-05"

$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Ly86"
$ws.Cells.Item(5,3).Value = "Cd180"
$ws.Cells.Item(5,4).Value = "Neutrophils"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(5,7).Value = [double]"0.641009"
$ws.Cells.Item(5,8).Value = [double]"1.923027"
$ws.Cells.Item(5,9).Value = [double]"0.005069072203501787"
$ws.Cells.Item(5,10).Value = [double]"0.005069072203501787"
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = [double]"43.85762"
$ws.Cells.Item(5,14).Value = [double]"131.57286"
$ws.Cells.Item(5,15).Value = [double]"0.1647163138771006"
$ws.Cells.Item(5,16).Value = [double]"0.1647163138771006"
$ws.Cells.Item(5,17).Value = [double]"28.11312913858"
$ws.Cells.Item(5,18).Value = [double]"253.01816224722"
$ws.Cells.Item(5,19).Value = [double]"0.0008349588881376863"
$ws.Cells.Item(5,20).Value = [double]"0.0008349588881376863"

$ws.Cells.Item(6,1).Value = "ECs"
$ws.Cells.Item(6,2).Value = "Ly86"
$ws.Cells.Item(6,3).Value = "Cd180"
$ws.Cells.Item(6,4).Value = "Resolving-Mac"
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(6,7).Value = [double]"0.641009"
$ws.Cells.Item(6,8).Value = [double]"1.923027"
$ws.Cells.Item(6,9).Value = [double]"0.005069072203501787"
$ws.Cells.Item(6,10).Value = [double]"0.005069072203501787"
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = [double]"146.7210513333333"
$ws.Cells.Item(6,14).Value = [double]"440.163154"
$ws.Cells.Item(6,15).Value = [double]"0.5510410903236319"
$ws.Cells.Item(6,16).Value = [double]"0.5510410903236319"
$ws.Cells.Item(6,17).Value = [double]"94.04951439412868"
$ws.Cells.Item(6,18).Value = [double]"846.4456295471581"
$ws.Cells.Item(6,19).Value = [double]"0.00279326707394684"
$ws.Cells.Item(6,20).Value = [double]"0.00279326707394684"

$ws.Cells.Item(7,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(7,2).Value = "Ly86"
$ws.Cells.Item(7,3).Value = "Cd180"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = [double]"65.36971666666666"
$ws.Cells.Item(7,8).Value = [double]"196.10915"
$ws.Cells.Item(7,9).Value = [double]"0.5169409691685881"
$ws.Cells.Item(7,10).Value = [double]"0.5169409691685881"
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = [double]"3.413185"
$ws.Cells.Item(7,14).Value = [double]"10.239555"
$ws.Cells.Item(7,15).Value = [double]"0.01281891839503857"
$ws.Cells.Item(7,16).Value = [double]"0.01281891839503857"
$ws.Cells.Item(7,17).Value = [double]"223.1189363809166"
$ws.Cells.Item(7,18).Value = [double]"2008.07042742825"
$ws.Cells.Item(7,19).Value = [double]"0.006626624098824282"
$ws.Cells.Item(7,20).Value = [double]"0.006626624098824282"

$ws.Cells.Item(8,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(8,2).Value = "Ly86"
$ws.Cells.Item(8,3).Value = "Cd180"
$ws.Cells.Item(8,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = [double]"65.36971666666666"
$ws.Cells.Item(8,8).Value = [double]"196.10915"
$ws.Cells.Item(8,9).Value = [double]"0.5169409691685881"
$ws.Cells.Item(8,10).Value = [double]"0.5169409691685881"
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = [double]"71.73578866666666"
$ws.Cells.Item(8,14).Value = [double]"215.207366"
$ws.Cells.Item(8,15).Value = [double]"0.2694185111330716"
$ws.Cells.Item(8,16).Value = [double]"0.2694185111330716"
$ws.Cells.Item(8,17).Value = [double]"4689.348179999877"
$ws.Cells.Item(8,18).Value = [double]"42204.13361999889"
$ws.Cells.Item(8,19).Value = [double]"0.1392734662570881"
$ws.Cells.Item(8,20).Value = [double]"0.1392734662570881"

$ws.Cells.Item(9,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(9,2).Value = "Ly86"
$ws.Cells.Item(9,3).Value = "Cd180"
$ws.Cells.Item(9,4).Value = "MuSCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = [double]"65.36971666666666"
$ws.Cells.Item(9,8).Value = [double]"196.10915"
$ws.Cells.Item(9,9).Value = [double]"0.5169409691685881"
$ws.Cells.Item(9,10).Value = [double]"0.5169409691685881"
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = [double]"0.5338986666666666"
$ws.Cells.Item(9,14).Value = [double]"1.601696"
$ws.Cells.Item(9,15).Value = [double]"0.002005166271157263"
$ws.Cells.Item(9,16).Value = [double]"0.002005166271157263"
$ws.Cells.Item(9,17).Value = [double]"34.90080456871111"
$ws.Cells.Item(9,18).Value = [double]"314.1072411184"
$ws.Cells.Item(9,19).Value = [double]"0.001036552595556199"
$ws.Cells.Item(9,20).Value = [double]"0.001036552595556199"

$ws.Cells.Item(10,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10,2).Value = "Ly86"
$ws.Cells.Item(10,3).Value = "Cd180"
$ws.Cells.Item(10,4).Value = "Neutrophils"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = [double]"65.36971666666666"
$ws.Cells.Item(10,8).Value = [double]"196.10915"
$ws.Cells.Item(10,9).Value = [double]"0.5169409691685881"
$ws.Cells.Item(10,10).Value = [double]"0.5169409691685881"
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = [double]"43.85762"
$ws.Cells.Item(10,14).Value = [double]"131.57286"
$ws.Cells.Item(10,15).Value = [double]"0.1647163138771006"
$ws.Cells.Item(10,16).Value = [double]"0.1647163138771006"
$ws.Cells.Item(10,17).Value = [double]"2866.960193074333"
$ws.Cells.Item(10,18).Value = [double]"25802.641737669"
$ws.Cells.Item(10,19).Value = [double]"0.08514861093350574"
$ws.Cells.Item(10,20).Value = [double]"0.08514861093350576"

$ws.Cells.Item(11,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(11,2).Value = "Ly86"
$ws.Cells.Item(11,3).Value = "Cd180"
$ws.Cells.Item(11,4).Value = "Resolving-Mac"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = [double]"65.36971666666666"
$ws.Cells.Item(11,8).Value = [double]"196.10915"
$ws.Cells.Item(11,9).Value = [double]"0.5169409691685881"
$ws.Cells.Item(11,10).Value = [double]"0.5169409691685881"
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = [double]"146.7210513333333"
$ws.Cells.Item(11,14).Value = [double]"440.163154"
$ws.Cells.Item(11,15).Value = [double]"0.5510410903236319"
$ws.Cells.Item(11,16).Value = [double]"0.5510410903236319"
$ws.Cells.Item(11,17).Value = [double]"9591.113554695456"
$ws.Cells.Item(11,18).Value = [double]"86320.0219922591"
$ws.Cells.Item(11,19).Value = [double]"0.2848557152836138"
$ws.Cells.Item(11,20).Value = [double]"0.2848557152836138"

$ws.Cells.Item(12,1).Value = "Neutrophils"
$ws.Cells.Item(12,2).Value = "Ly86"
$ws.Cells.Item(12,3).Value = "Cd180"
$ws.Cells.Item(12,4).Value = "ECs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = [double]"10.18446633333333"
$ws.Cells.Item(12,8).Value = [double]"30.553399"
$ws.Cells.Item(12,9).Value = [double]"0.08053833128364775"
$ws.Cells.Item(12,10).Value = [double]"0.08053833128364775"
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = [double]"3.413185"
$ws.Cells.Item(12,14).Value = [double]"10.239555"
$ws.Cells.Item(12,15).Value = [double]"0.01281891839503857"
$ws.Cells.Item(12,16).Value = [double]"0.01281891839503857"
$ws.Cells.Item(12,17).Value = [double]"34.76146772193833"
$ws.Cells.Item(12,18).Value = [double]"312.853209497445"
$ws.Cells.Item(12,19).Value = [double]"0.001032414296397663"
$ws.Cells.Item(12,20).Value = [double]"0.001032414296397663"

$ws.Cells.Item(13,1).Value = "Neutrophils"
$ws.Cells.Item(13,2).Value = "Ly86"
$ws.Cells.Item(13,3).Value = "Cd180"
$ws.Cells.Item(13,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = [double]"10.18446633333333"
$ws.Cells.Item(13,8).Value = [double]"30.553399"
$ws.Cells.Item(13,9).Value = [double]"0.08053833128364775"
$ws.Cells.Item(13,10).Value = [double]"0.08053833128364775"
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = [double]"71.73578866666666"
$ws.Cells.Item(13,14).Value = [double]"215.207366"
$ws.Cells.Item(13,15).Value = [double]"0.2694185111330716"
$ws.Cells.Item(13,16).Value = [double]"0.2694185111330716"
$ws.Cells.Item(13,17).Value = [double]"730.5907245707815"
$ws.Cells.Item(13,18).Value = [double]"6575.316521137033"
$ws.Cells.Item(13,19).Value = [double]"0.02169851730358246"
$ws.Cells.Item(13,20).Value = [double]"0.02169851730358246"

$ws.Cells.Item(14,1).Value = "Neutrophils"
$ws.Cells.Item(14,2).Value = "Ly86"
$ws.Cells.Item(14,3).Value = "Cd180"
$ws.Cells.Item(14,4).Value = "MuSCs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = [double]"10.18446633333333"
$ws.Cells.Item(14,8).Value = [double]"30.553399"
$ws.Cells.Item(14,9).Value = [double]"0.08053833128364775"
$ws.Cells.Item(14,10).Value = [double]"0.08053833128364775"
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = [double]"0.5338986666666666"
$ws.Cells.Item(14,14).Value = [double]"1.601696"
$ws.Cells.Item(14,15).Value = [double]"0.002005166271157263"
$ws.Cells.Item(14,16).Value = [double]"0.002005166271157263"
$ws.Cells.Item(14,17).Value = [double]"5.437472996078221"
$ws.Cells.Item(14,18).Value = [double]"48.937256964704"
$ws.Cells.Item(14,19).Value = [double]"0.0001614927454252603"
$ws.Cells.Item(14,20).Value = [double]"0.0001614927454252603"

$ws.Cells.Item(15,1).Value = "Neutrophils"
$ws.Cells.Item(15,2).Value = "Ly86"
$ws.Cells.Item(15,3).Value = "Cd180"
$ws.Cells.Item(15,4).Value = "Neutrophils"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = [double]"10.18446633333333"
$ws.Cells.Item(15,8).Value = [double]"30.553399"
$ws.Cells.Item(15,9).Value = [double]"0.08053833128364775"
$ws.Cells.Item(15,10).Value = [double]"0.08053833128364775"
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = [double]"43.85762"
$ws.Cells.Item(15,14).Value = [double]"131.57286"
$ws.Cells.Item(15,15).Value = [double]"0.1647163138771006"
$ws.Cells.Item(15,16).Value = [double]"0.1647163138771006"
$ws.Cells.Item(15,17).Value = [double]"446.6664543501266"
$ws.Cells.Item(15,18).Value = [double]"4019.99808915114"
$ws.Cells.Item(15,19).Value = [double]"0.01326597705485523"
$ws.Cells.Item(15,20).Value = [double]"0.01326597705485523"

$ws.Cells.Item(16,1).Value = "Neutrophils"
$ws.Cells.Item(16,2).Value = "Ly86"
$ws.Cells.Item(16,3).Value = "Cd180"
$ws.Cells.Item(16,4).Value = "Resolving-Mac"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = [double]"10.18446633333333"
$ws.Cells.Item(16,8).Value = [double]"30.553399"
$ws.Cells.Item(16,9).Value = [double]"0.08053833128364775"
$ws.Cells.Item(16,10).Value = [double]"0.08053833128364775"
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = [double]"146.7210513333333"
$ws.Cells.Item(16,14).Value = [double]"440.163154"
$ws.Cells.Item(16,15).Value = [double]"0.5510410903236319"
$ws.Cells.Item(16,16).Value = [double]"0.5510410903236319"
$ws.Cells.Item(16,17).Value = [double]"1494.275607695605"
$ws.Cells.Item(16,18).Value = [double]"13448.48046926045"
$ws.Cells.Item(16,19).Value = [double]"0.04437992988338713"
$ws.Cells.Item(16,20).Value = [double]"0.04437992988338713"

$ws.Cells.Item(17,1).Value = "Resolving-Mac"
$ws.Cells.Item(17,2).Value = "Ly86"
$ws.Cells.Item(17,3).Value = "Cd180"
$ws.Cells.Item(17,4).Value = "ECs"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = [double]"50.25970433333333"
$ws.Cells.Item(17,8).Value = [double]"150.779113"
$ws.Cells.Item(17,9).Value = [double]"0.3974516273442624"
$ws.Cells.Item(17,10).Value = [double]"0.3974516273442624"
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = [double]"3.413185"
$ws.Cells.Item(17,14).Value = [double]"10.239555"
$ws.Cells.Item(17,15).Value = [double]"0.01281891839503857"
$ws.Cells.Item(17,16).Value = [double]"0.01281891839503857"
$ws.Cells.Item(17,17).Value = [double]"171.5456689349683"
$ws.Cells.Item(17,18).Value = [double]"1543.911020414715"
$ws.Cells.Item(17,19).Value = [double]"0.005094899976901382"
$ws.Cells.Item(17,20).Value = [double]"0.005094899976901381"

$ws.Cells.Item(18,1).Value = "Resolving-Mac"
$ws.Cells.Item(18,2).Value = "Ly86"
$ws.Cells.Item(18,3).Value = "Cd180"
$ws.Cells.Item(18,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(18,5).Value = 3
$ws.Cells.Item(18,6).Value = 1
$ws.Cells.Item(18,7).Value = [double]"50.25970433333333"
$ws.Cells.Item(18,8).Value = [double]"150.779113"
$ws.Cells.Item(18,9).Value = [double]"0.3974516273442624"
$ws.Cells.Item(18,10).Value = [double]"0.3974516273442624"
$ws.Cells.Item(18,11).Value = 3
$ws.Cells.Item(18,12).Value = 1
$ws.Cells.Item(18,13).Value = [double]"71.73578866666666"
$ws.Cells.Item(18,14).Value = [double]"215.207366"
$ws.Cells.Item(18,15).Value = [double]"0.2694185111330716"
$ws.Cells.Item(18,16).Value = [double]"0.2694185111330716"
$ws.Cells.Item(18,17).Value = [double]"3605.419528505151"
$ws.Cells.Item(18,18).Value = [double]"32448.77575654635"
$ws.Cells.Item(18,19).Value = [double]"0.1070808256865076"
$ws.Cells.Item(18,20).Value = [double]"0.1070808256865076"

$ws.Cells.Item(19,1).Value = "Resolving-Mac"
$ws.Cells.Item(19,2).Value = "Ly86"
$ws.Cells.Item(19,3).Value = "Cd180"
$ws.Cells.Item(19,4).Value = "MuSCs"
$ws.Cells.Item(19,5).Value = 3
$ws.Cells.Item(19,6).Value = 1
$ws.Cells.Item(19,7).Value = [double]"50.25970433333333"
$ws.Cells.Item(19,8).Value = [double]"150.779113"
$ws.Cells.Item(19,9).Value = [double]"0.3974516273442624"
$ws.Cells.Item(19,10).Value = [double]"0.3974516273442624"
$ws.Cells.Item(19,11).Value = 3
$ws.Cells.Item(19,12).Value = 1
$ws.Cells.Item(19,13).Value = [double]"0.5338986666666666"
$ws.Cells.Item(19,14).Value = [double]"1.601696"
$ws.Cells.Item(19,15).Value = [double]"0.002005166271157263"
$ws.Cells.Item(19,16).Value = [double]"0.002005166271157263"
$ws.Cells.Item(19,17).Value = [double]"26.83358913062755"
$ws.Cells.Item(19,18).Value = [double]"241.502302175648"
$ws.Cells.Item(19,19).Value = [double]"0.0007969565975672805"
$ws.Cells.Item(19,20).Value = [double]"0.0007969565975672804"

$ws.Cells.Item(20,1).Value = "Resolving-Mac"
$ws.Cells.Item(20,2).Value = "Ly86"
$ws.Cells.Item(20,3).Value = "Cd180"
$ws.Cells.Item(20,4).Value = "Neutrophils"
$ws.Cells.Item(20,5).Value = 3
$ws.Cells.Item(20,6).Value = 1
$ws.Cells.Item(20,7).Value = [double]"50.25970433333333"
$ws.Cells.Item(20,8).Value = [double]"150.779113"
$ws.Cells.Item(20,9).Value = [double]"0.3974516273442624"
$ws.Cells.Item(20,10).Value = [double]"0.3974516273442624"
$ws.Cells.Item(20,11).Value = 3
$ws.Cells.Item(20,12).Value = 1
$ws.Cells.Item(20,13).Value = [double]"43.85762"
$ws.Cells.Item(20,14).Value = [double]"131.57286"
$ws.Cells.Item(20,15).Value = [double]"0.1647163138771006"
$ws.Cells.Item(20,16).Value = [double]"0.1647163138771006"
$ws.Cells.Item(20,17).Value = [double]"2204.271013963687"
$ws.Cells.Item(20,18).Value = [double]"19838.43912567318"
$ws.Cells.Item(20,19).Value = [double]"0.06546676700060194"
$ws.Cells.Item(20,20).Value = [double]"0.06546676700060194"

$ws.Cells.Item(21,1).Value = "Resolving-Mac"
$ws.Cells.Item(21,2).Value = "Ly86"
$ws.Cells.Item(21,3).Value = "Cd180"
$ws.Cells.Item(21,4).Value = "Resolving-Mac"
$ws.Cells.Item(21,5).Value = 3
$ws.Cells.Item(21,6).Value = 1
$ws.Cells.Item(21,7).Value = [double]"50.25970433333333"
$ws.Cells.Item(21,8).Value = [double]"150.779113"
$ws.Cells.Item(21,9).Value = [double]"0.3974516273442624"
$ws.Cells.Item(21,10).Value = [double]"0.3974516273442624"
$ws.Cells.Item(21,11).Value = 3
$ws.Cells.Item(21,12).Value = 1
$ws.Cells.Item(21,13).Value = [double]"146.7210513333333"
$ws.Cells.Item(21,14).Value = [double]"440.163154"
$ws.Cells.Item(21,15).Value = [double]"0.5510410903236319"
$ws.Cells.Item(21,16).Value = [double]"0.5510410903236319"
$ws.Cells.Item(21,17).Value = [double]"7374.156659489156"
$ws.Cells.Item(21,18).Value = [double]"66367.40993540241"
$ws.Cells.Item(21,19).Value = [double]"0.2190121780826842"
$ws.Cells.Item(21,20).Value = [double]"0.2190121780826842"
